# "versao final final com 90 MHZ"
# Update the two fuzzy-logic controller inputs on Planilha1 (A3, B3).
# All downstream formulas (E3:P3, A8:C8, A11:C11, F9, ...) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("A3").Value = 160
$ws.Range("B3").Value = 80
